$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item(1)

# Add a new column "Link of Blog Post on The Hacked Site" to the table (col G)
$tbl.ListColumns.Add() | Out-Null

# Add a new row to the table (row 94) for Post 84
$newRow = $tbl.ListRows.Add()

# Carry formatting from the previous data row (93) into the new row
$ws.Range("B93").Copy($ws.Range("B94"))
$ws.Range("C93").Copy($ws.Range("C94"))
$ws.Range("D93").Copy($ws.Range("D94"))
$ws.Range("E93").Copy($ws.Range("E94"))
$ws.Range("F93").Copy($ws.Range("F94"))

# Fill in the data for Post 84 - "File Encrypter and Decrypter"
$ws.Range("B94").Value = 84
$ws.Range("C94").Value = "File Encrypter and Decrypter "
$ws.Range("D94").Value = 44263
$ws.Range("F94").Value = "https://dev.to/rahulmishra05/file-encrypter-and-decrypter-10mo"
$ws.Range("E94").Value = "https://programmingport.hashnode.dev/file-encrypter-and-decrypter"
$ws.Range("G94").Value = "https://thehackedsite.netlify.app/shell/script/2021/03/08/file-encrypter-decrypter"

# Header for the new column
$ws.Range("G10").Value = "Link of Blog Post on The Hacked Site"

# Column width for the new column G
$ws.Columns.Item(7).ColumnWidth = 85.65

# Update the active selection to the newly added cell
$ws.Range("G94").Select() | Out-Null
